$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of B2:E51 are plain text cells in the source data (coin name, link,
# price, % change). Several "Price" values are digit strings that Excel's
# General-format parser would otherwise coerce into numbers (losing e.g. a
# trailing zero, "0.07750" -> 0.0775), so those are entered with a leading
# apostrophe to force literal text, matching the original text content.

$ws.Range("D2").Value = "25.861.12"
$ws.Range("D3").Value = "1.639.29"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -1.08%  "
$ws.Range("D5").Value = "'215.43"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'0.5033"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").Value = "'0.2575"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").Value = "'0.06376"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "'19.49"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").Value = "'0.07750"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.654.58"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.259"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.867.19"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "'0.5453"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "0.0₅7891"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "'64.19"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "25.905.91"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").Value = "'202.11"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").Value = "'4.386"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'9.883"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").Value = "'5.970"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "'1.864"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "'140.92"
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").Value = "'0.1135"
$ws.Range("E27").Value = "  -3.15%  "
$ws.Range("D28").Value = "'15.67"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'6.777"
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("D30").Value = "'1.244"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "'0.04972"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("D32").Value = "'3.271"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").Value = "'3.195"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "'1.547"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "'2.370"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").Value = "'2.627"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").Value = "'0.8909"
$ws.Range("E37").Value = "  -3.18%  "
$ws.Range("D38").Value = "1.149.49"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").Value = "'0.5593"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("D40").Value = "'0.01564"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "'1.004"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").Value = "'5.685"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'99.77"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8059"
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").Value = "1.778.17"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("D47").Value = "'0.4525"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").Value = "'54.74"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").Value = "'0.05049"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  -0.87%  "
